$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$updates = [ordered]@{
    "D2" = "30.941.80"
    "E2" = "  +2.73%  "
    "D3" = "2.120.95"
    "E3" = "  +10.61%  "
    "D4" = "1.001"
    "E4" = "  -0.26%  "
    "D5" = "334.55"
    "E5" = "  +4.80%  "
    "D6" = "0.9999"
    "E6" = "  -0.21%  "
    "D7" = "0.5355"
    "E7" = "  +5.68%  "
    "D8" = "0.4417"
    "E8" = "  +8.40%  "
    "D9" = "0.09077"
    "E9" = "  +8.93%  "
    "D10" = "46.45"
    "E10" = "  +10.51%  "
    "D11" = "1.184"
    "E11" = "  +6.00%  "
    "D12" = "25.43"
    "E12" = "  +4.83%  "
    "D13" = "2.118.65"
    "E13" = "  +10.28%  "
    "D14" = "6.804"
    "E14" = "  +5.94%  "
    "D15" = "7.841"
    "E15" = "  +8.14%  "
    "D16" = "98.25"
    "E16" = "  +6.06%  "
    "B17" = "ShibaInu"
    "C17" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D17" = "0.00001140"
    "E17" = "  +4.11%  "
    "B18" = "BinanceUSD"
    "C18" = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
    "D18" = "1.001"
    "E18" = "  -0.32%  "
    "D19" = "0.06666"
    "E19" = "  +2.29%  "
    "D20" = "19.27"
    "E20" = "  +4.02%  "
    "D21" = "0.9993"
    "E21" = "  -0.21%  "
    "D22" = "6.404"
    "E22" = "  +7.58%  "
    "D23" = "31.020.61"
    "E23" = "  +2.96%  "
    "E24" = "  +7.46%  "
    "D25" = "2.368.27"
    "D26" = "2.269"
    "E26" = "  +3.30%  "
    "D27" = "22.94"
    "E27" = "  +4.99%  "
    "E28" = "  +14.71%  "
    "D29" = "163.65"
    "E29" = "  +0.52%  "
    "D30" = "134.30"
    "E30" = "  +4.31%  "
    "D31" = "1.181"
    "E31" = "  +3.88%  "
    "D32" = "0.1081"
    "E32" = "  +3.26%  "
    "D33" = "6.279"
    "E33" = "  +5.59%  "
    "D34" = "4.011"
    "E34" = "  +5.78%  "
    "D35" = "1.540"
    "E35" = "  +27.27%  "
    "D36" = "0.02618"
    "E36" = "  +6.92%  "
    "D37" = "13.34"
    "E37" = "  +16.33%  "
    "D38" = "5.595"
    "E38" = "  +5.67%  "
    "D39" = "9.624"
    "E39" = "  +12.26%  "
    "D40" = "0.06775"
    "E40" = "  +5.33%  "
    "D41" = "0.2280"
    "E41" = "  +6.41%  "
    "D42" = "0.6882"
    "E42" = "  +6.20%  "
    "E43" = "  +3.96%  "
    "B44" = "Decentraland"
    "C44" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D44" = "0.6480"
    "E44" = "  +7.25%  "
    "B45" = "EnergySwap"
    "C45" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D45" = "14.18"
    "E45" = "  +6.18%  "
    "D46" = "0.9993"
    "E46" = "  -0.10%  "
    "D47" = "2.263"
    "E47" = "  +3.56%  "
    "D48" = "3.681"
    "E48" = "  +1.58%  "
    "E49" = "  +6.43%  "
    "D50" = "83.54"
    "E50" = "  +7.36%  "
    "D51" = "1.182"
    "E51" = "  +4.22%  "
}

# Cells in D (Price) and E (Volume) columns must stay text even when the
# new value looks like a pure number (Excel would otherwise auto-convert it).
foreach ($addr in $updates.Keys) {
    $col = $addr.Substring(0,1)
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($addr).NumberFormat = "@"
    }
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore default (Normal) cell style so no stray number-format styles remain
foreach ($addr in $updates.Keys) {
    $col = $addr.Substring(0,1)
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($addr).Style = "Normal"
    }
}
